$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22, column A currently stores the phone number "79174445" as text.
# Re-enter it as a genuine number (same value), matching the rest of the
# "phone" column.
$ws.Range("A22").Value = 79174445

# Append the new redemption (row 23): phone, points redeemed, timestamp.
# Phone numbers are stored as text in this sheet, so force column A back
# to Text before writing the numeric-looking string, otherwise Excel will
# auto-coerce it to a number like it just did above for A22.
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "79174445"
$ws.Range("A23").NumberFormat = "General"
$ws.Range("A23").Style = "Normal"

$ws.Range("B23").Value = 500
$ws.Range("C23").Value = "2025-08-18T17:43:28"
